# Applies the LinuxForHealth re-brand / version-bump edit described by the
# commit "Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig@80fa500..."
#
# Resolved cell-level changes:
#   Metadata!B2  URL           http://ibm.com/...           -> http://linuxforhealth.org/...
#   Metadata!B3  Version       7.0.0                         -> 8.0.0
#   Metadata!B8  Date          2022-09-08T16:11:15+00:00     -> 2022-11-10T16:00:46+00:00
#   Metadata!B9  Publisher     Alvearie Team                 -> LinuxForHealth Team
#   Elements!Q5  Fixed Value   http://ibm.com/...            -> http://linuxforhealth.org/... (same text as B2)
#   Elements!AI2 Constraint(s) "ele-1:...\next-1:..."        -> "" (empty text; this note now only
#                 appears on the Extension.extension row, AI4, which already carries it unchanged)

$wb     = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsEl   = $wb.Worksheets.Item("Elements")

$oldUrl = "http://ibm.com/fhir/cdm/StructureDefinition/detected"
$newUrl = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/detected"

# Update the URL wherever it appears (Metadata!B2 and Elements!Q5 share the
# same text) using Replace so both cells stay backed by one shared string,
# matching how the workbook stores it.
$wsMeta.Cells.Replace($oldUrl, $newUrl) | Out-Null
$wsEl.Cells.Replace($oldUrl, $newUrl) | Out-Null

# Version bump
$wsMeta.Range("B3").Value = "8.0.0"

# Date bump
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher rename
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# The "ele-1 / ext-1" constraint note on the root Extension row (AI2) is
# cleared -- in the updated IG it is only listed against Extension.extension
# (AI4), which already holds this text and is left untouched. A plain
# Value = "" assignment turns the cell fully blank (numeric/empty) instead
# of an empty *text* cell, so we prefix with a single quote (Excel's
# text-literal marker, stripped from the stored value) to force text type,
# then re-paste the original cell's formatting so no stray number format /
# quote-prefix style sticks to the cell.
$aiCell = $wsEl.Range("AI2")
$aiCell.Value = "'"
$wsEl.Range("AH2").Copy()
$aiCell.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
